$d = $word.ActiveDocument

# The target edit only adds a built-in "Hyperlink" character style definition
# to styles.xml (no document body changes). We mint that style definition by
# momentarily applying the "Hyperlink" style to a throwaway run of text in a
# scratch paragraph, then remove the scratch paragraph again so the visible
# document content is left untouched.
$scratch = $d.Range(0, 0)
$scratch.InsertParagraphBefore()
$scratchPara = $d.Paragraphs(1)
$scratchPara.Range.Text = "x"
$scratchRun = $d.Range($scratchPara.Range.Start, $scratchPara.Range.Start + 1)
$scratchRun.Style = "Hyperlink"
$scratchPara.Range.Delete()

# Bring the newly minted style's formatting in line with the classic
# built-in "Hyperlink" character style (blue, single-underline, uiPriority 99,
# semi-hidden until used).
$hyperlinkStyle = $d.Styles("Hyperlink")
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true
$hyperlinkStyle.QuickStyle = $false
$hyperlinkStyle.Font.Color = 255 * 65536
$hyperlinkStyle.Font.Underline = 1
